$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.778.31"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "3.407.38"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.60"
$ws.Range("D6").ClearFormats()

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.408.59"
$ws.Range("E8").Value = "  -0.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.548"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -9.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("E11").Value = "  -3.56%  "

$ws.Range("E12").Value = "  -4.46%  "

$ws.Range("D13").Value = "3.996.68"
$ws.Range("E13").Value = "  -0.37%  "


$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.98"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.16%  "

$ws.Range("E16").Value = "  -9.04%  "

$ws.Range("D17").Value = "63.814.91"
$ws.Range("E17").Value = "  -1.33%  "

$ws.Range("D18").Value = "3.408.69"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.05"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.49"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.25"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.35%  "

$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.39"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.514"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.30%  "

$ws.Range("E26").Value = "  -2.84%  "

$ws.Range("E27").Value = "  -5.30%  "

$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  -3.82%  "

$ws.Range("E31").Value = "  -6.11%  "

$ws.Range("E32").Value = "  -0.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.76"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.01"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.51"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.90"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.86%  "

$ws.Range("E37").Value = "  +7.99%  "

$ws.Range("E38").Value = "  -6.32%  "

$ws.Range("D39").Value = "2.809.47"
$ws.Range("E39").Value = "  -2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0724"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.98"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.50"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.61"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.03"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.39"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.36%  "

$ws.Range("E46").Value = "  -4.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "335.22"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.38%  "

$ws.Range("E49").Value = "  -2.43%  "

$ws.Range("E50").Value = "  -4.74%  "

$ws.Range("E51").Value = "  -4.26%  "
